$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New LocationY values (column B, rows 2-62) after applying the 3D
# transformation pipeline and filtering out z coordinates below 0.
$newLocationY = @(
    -20.09988556775736,
    26.93396066402152,
    26.1890690863497,
    40.5294689487978,
    57.45701136205512,
    68.08166002029758,
    92.38384244685979,
    101.8395708829993,
    118.3588582572152,
    127.9501875438875,
    135.1299865423651,
    147.1257779744952,
    163.1948116121281,
    171.9127730373955,
    183.8433836327171,
    191.7349824875643,
    202.7406588936277,
    208.4263569882514,
    220.955472444577,
    226.6929536597885,
    234.8831737574142,
    222.997923082522,
    253.0106256471292,
    265.0575906020436,
    271.485941374627,
    285.7487970790505,
    301.9355845705513,
    304.9168434895126,
    316.1300149582157,
    331.0114815109552,
    331.1661896002388,
    342.2958775762802,
    387.1280950623595,
    402.6338330601983,
    428.7572957296376,
    440.4314262612564,
    471.5973450697712,
    478.8135481218333,
    494.5974419764873,
    506.4176502373559,
    519.8145692529029,
    523.1197761242677,
    515.4281346601889,
    543.6370417721078,
    552.8921870723887,
    566.5167712529859,
    571.0145799563303,
    599.6063612714547,
    606.504222669016,
    613.9769481991199,
    633.5287141327602,
    639.9936875478271,
    648.8585111159709,
    663.3866026351061,
    659.2658763869146,
    676.1422081189438,
    677.5446386290897,
    675.8333807842894,
    685.0206265226069,
    684.7600759324165,
    685.8514501910832
)

for ($i = 0; $i -lt $newLocationY.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value2 = $newLocationY[$i]
}
